$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(9)
$tr = $sh.TextFrame.TextRange
$paras = $tr.Paragraphs()
for ($i = 1; $i -le $paras.Count; $i++) {
    $para = $paras.Item($i)
    Write-Output ("{0}: [{1}] indentLevel={2}" -f $i, $para.Text, $para.IndentLevel)
}
